# Daily attendance processing - 2025-10-22 07:21:07
#
# Column G ("Recorded By") lists the users who touched each attendance
# session, e.g. "System, dnasr281@gmail.com". This pass normalises the
# ordering of that list by reversing it (most-recent-actor-first style),
# e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System".
#
# Only cells holding more than one comma-separated name are affected;
# single-name cells, blank cells, and any entry that still references the
# admin@admin.com placeholder account are left exactly as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in the "Recorded By" column (G) that contain more than one
# comma-separated entry and need that entry order reversed.
$rowsToReverse = @(
    2, 3, 4, 5, 6, 8, 10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22,
    29, 30, 31, 32, 33, 35, 37, 38, 39, 40, 41, 42, 44, 45, 46, 47, 48, 49,
    56, 57, 58, 59, 60, 62, 64, 65, 66, 67, 68, 69, 71, 72, 73, 74, 75, 76,
    83, 84, 85, 86, 87, 88, 89, 93, 95, 96, 97, 99,
    109, 110, 111, 112, 113, 114, 115, 119, 121, 122, 123, 125,
    135, 136, 137, 138, 139, 140, 141, 145, 147, 148, 149, 151
)

foreach ($r in $rowsToReverse) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",") -and -not $val.Contains("admin@admin.com")) {
        $parts = $val -split ",\s*"
        if ($parts.Count -gt 1) {
            $reversed = $parts[($parts.Count - 1)..0]
            $cell.Value2 = [string]::Join(", ", $reversed)
        }
    }
}
